# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46055
$ws.Range("B2").Value = -0.07000000000000001
$ws.Range("C2").Value = -0.16
$ws.Range("D2").Value = -0.23
$ws.Range("E2").Value = -0.5
$ws.Range("F2").Value = -0.5
$ws.Range("G2").Value = -0.5
$ws.Range("H2").Value = -0.5
$ws.Range("I2").Value = -0.17
$ws.Range("J2").Value = -0.01
$ws.Range("K2").Value = 0.06
$ws.Range("L2").Value = 0.22
$ws.Range("M2").Value = 0.06
$ws.Range("N2").Value = 0.08
$ws.Range("O2").Value = 0.03
$ws.Range("P2").Value = 0.01
$ws.Range("Q2").Value = 0.02
$ws.Range("R2").Value = 0.88
$ws.Range("S2").Value = 0.4
$ws.Range("T2").Value = 5.15
$ws.Range("U2").Value = 15.7
$ws.Range("V2").Value = 28.97
$ws.Range("W2").Value = 28.97
$ws.Range("X2").Value = 22.37
$ws.Range("Y2").Value = 6.66
$ws.Range("Z2").Value = 4.46
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 21.74
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 28.97
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 14.52
$ws.Range("AG2").Value = "0h-17h"
